# Update "想去人数" (interest count) figures in the 展览 (Exhibition) and
# 全部类型 (All types) sheets to match the latest generated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 7049
$ws1.Range("F4").Value  = 71
$ws1.Range("F5").Value  = 460
$ws1.Range("F7").Value  = 6994
$ws1.Range("F11").Value = 25
$ws1.Range("F14").Value = 154
$ws1.Range("F16").Value = 418
$ws1.Range("F18").Value = 50
$ws1.Range("F19").Value = 21
$ws1.Range("F20").Value = 5354
$ws1.Range("F22").Value = 187
$ws1.Range("F23").Value = 778

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 7049
$ws4.Range("F4").Value  = 71
$ws4.Range("F5").Value  = 460
$ws4.Range("F7").Value  = 6994
$ws4.Range("F11").Value = 25
$ws4.Range("F14").Value = 154
$ws4.Range("F16").Value = 418
$ws4.Range("F18").Value = 50
$ws4.Range("F19").Value = 21
$ws4.Range("F21").Value = 5354
$ws4.Range("F24").Value = 187
$ws4.Range("F25").Value = 778
